$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename test fixture from "ScopeProperties*" to "HeadersToClaims*" / "Claims*" ---
# Column A (ProjectName): ScopePropertiesApi -> HeadersToClaimsApi
$ws.Range("A2:A7").Value = "HeadersToClaimsApi"

# Column B (ClassName): ScopePropertiesController -> ClaimsController
$ws.Range("B2:B7").Value = "ClaimsController"

# --- Update expected JSON claim payloads: claim Types changed from header names
#     (X-UserScope / X-Role / X-User) to claim names (user_scope / role / name),
#     and the combined row now orders name/role before user_scope ---
$ws.Range("G3").Value = '[{"Type":"user_scope","Value":"ABC"}]'
$ws.Range("G5").Value = '[{"Type":"name","Value":"moe@stooges.org"},{"Type":"role","Value":"admin"}]'
$ws.Range("G7").Value = '[{"Type":"name","Value":"moe@stooges.org"},{"Type":"role","Value":"admin"},{"Type":"user_scope","Value":"ABC"}]'

# --- Columns A:D got narrower to best-fit the new (shorter) text ---
$ws.Columns.Item(1).ColumnWidth = 16.333333333333336
$ws.Columns.Item(2).ColumnWidth = 13.166666666666668
$ws.Columns.Item(3).ColumnWidth = 10.833333333333332
$ws.Columns.Item(4).ColumnWidth = 10.0
